$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# --- Swap the "Arroyo de la Luz" / "La Gomera" rows (row 62 / row 63) ---
# Row 62 becomes "La Gomera" with updated figures
$ws.Range("A62").Value = "La Gomera"
$ws.Range("B62").Value = 9
$ws.Range("C62").Value = 4
$ws.Range("D62").Value = 5
$ws.Range("E62").Value = 0

# Row 63 becomes "Arroyo de la Luz"
$ws.Range("A63").Value = "Arroyo de la Luz"
$ws.Range("B63").Value = 7
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 7
$ws.Range("E63").Value = 0

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 22:22"

# --- Update Tenerife (row 32) figures ---
$ws.Range("B32").Value = 1140
$ws.Range("C32").Value = 216
$ws.Range("D32").Value = 857
$ws.Range("E32").Value = 67

# --- Update Gran Canaria (row 47) figures ---
$ws.Range("B47").Value = 434
$ws.Range("C47").Value = 104
$ws.Range("D47").Value = 305
$ws.Range("E47").Value = 25

# --- Update La Palma (row 56) figures ---
$ws.Range("B56").Value = 79
$ws.Range("C56").Value = 9
$ws.Range("D56").Value = 68
$ws.Range("E56").Value = 2

# --- Update Lanzarote (row 57) figures ---
$ws.Range("B57").Value = 72
$ws.Range("C57").Value = 14
$ws.Range("D57").Value = 55
$ws.Range("E57").Value = 3

# --- Update Fuerteventura (row 59) figures ---
$ws.Range("B59").Value = 40
$ws.Range("C59").Value = 16
$ws.Range("D59").Value = 24

# --- Update El Hierro (row 64) figures ---
$ws.Range("B64").Value = 4
$ws.Range("D64").Value = 3
